$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" (D) column stores plain formatted-price TEXT (not numbers) in the
# source data, matching the original workbook where every D-cell is an inline
# string. Assigning a plain numeric-looking string via .Value would make Excel
# auto-convert the cell to a Number, so for values that look numeric we briefly
# force Text number-formatting, assign the literal text, then restore the cell
# to the default "Normal" style so no stray formatting is left behind - only the
# text content changes, exactly like the other already-non-numeric price strings
# (e.g. "42.337.38") which Excel could never auto-convert in the first place.

$ws.Range("D2").Value = '42.337.38'
$ws.Range("E2").Value = '  -2.99%  '

$ws.Range("D3").Value = '2.209.54'
$ws.Range("E3").Value = '  -3.17%  '

$ws.Range("E4").Value = '  +0.57%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '107.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -12.90%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '295.23'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +10.83%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.13%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.602'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.73%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '44.39'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -8.49%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0915'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.38%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.29'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.11%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.76'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.28%  '

$ws.Range("E14").Value = '  -2.65%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.936'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.30%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.94'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.41%  '

$ws.Range("D17").Value = '2.551.25'

$ws.Range("D18").Value = '2.241.39'
$ws.Range("E18").Value = '  -1.75%  '

$ws.Range("D19").Value = '42.318.11'
$ws.Range("E19").Value = '  -3.27%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.23%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000105'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.71'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.78%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.41'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +18.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.28'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -7.20%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '228.75'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.22'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.57%  '

$ws.Range("E27").Value = '  -1.84%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.59'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.16%  '

$ws.Range("E29").Value = '  -1.00%  '

$ws.Range("E30").Value = '  -1.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '38.34'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -10.92%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.21'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '174.18'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.87%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '20.96'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.61%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0877'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.62'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.50%  '

$ws.Range("E39").Value = '  -3.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0364'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.88%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.102'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.95%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.46'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.83%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.233'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.21%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '70.34'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.29%  '

$ws.Range("E45").Value = '  +0.23%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.59'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -10.25%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.30'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.23%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.44'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.10%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.31'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.89%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.62'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.57%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.47'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.80%  '

# Rows 37/38: RenderToken and NEARProtocol swap rank positions
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.24'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.12%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.81'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.04%  '

